$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 and Row 5 have their observation-specific data swapped
# (columns A, I, J, Q, R, X, Z, AB). Everything else (B, D-H, P, S, T-W,
# Y, AA, AD-AY, ...) stays identical between the two rows.

# --- Capture row 4 "observation" values before overwriting ---
$A4 = $ws.Range("A4").Value2
$I4 = $ws.Range("I4").Value2
$J4 = $ws.Range("J4").Value2
$Q4 = $ws.Range("Q4").Value2
$R4 = $ws.Range("R4").Value2
$X4 = $ws.Range("X4").Value2
$Z4 = $ws.Range("Z4").Value2
$AB4 = $ws.Range("AB4").Value2

# --- Capture row 5 "observation" values before overwriting ---
$A5 = $ws.Range("A5").Value2
$I5 = $ws.Range("I5").Value2
$J5 = $ws.Range("J5").Value2
$Q5 = $ws.Range("Q5").Value2
$R5 = $ws.Range("R5").Value2
$X5 = $ws.Range("X5").Value2
$Z5 = $ws.Range("Z5").Value2
$AB5 = $ws.Range("AB5").Value2

# --- Write row 5's old values into row 4 ---
$ws.Range("A4").Value = $A5
$ws.Range("Q4").Value = $Q5
$ws.Range("R4").Value = $R5
$ws.Range("X4").Value = $X5
$ws.Range("Z4").Value = $Z5
$ws.Range("AB4").Value = $AB5

# I4/J4 become empty (row 5 had nothing in I5/J5)
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""

# --- Write row 4's old values into row 5 ---
$ws.Range("A5").Value = $A4
$ws.Range("Q5").Value = $Q4
$ws.Range("R5").Value = $R4
$ws.Range("X5").Value = $X4
$ws.Range("Z5").Value = $Z4
$ws.Range("AB5").Value = $AB4

# I5/J5 get row 4's old "10" / "bålar" (must stay text, not become numeric)
$ws.Range("I5").Value = "'" + $I4
$ws.Range("I5").Style = "Normal"
$ws.Range("J5").Value = $J4
